$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$styleD2 = $ws.Range("D2").Style
$ws.Range("D2").Value = '''59.903.13'
$ws.Range("D2").Style = $styleD2
$ws.Range("E2").Value = '  +3.34%  '
$styleD3 = $ws.Range("D3").Style
$ws.Range("D3").Value = '''2.341.07'
$ws.Range("D3").Style = $styleD3
$ws.Range("E3").Value = '  +2.46%  '
$ws.Range("E4").Value = '  -0.01%  '
$styleD5 = $ws.Range("D5").Style
$ws.Range("D5").Value = '''544.15'
$ws.Range("D5").Style = $styleD5
$ws.Range("E5").Value = '  +1.99%  '
$styleD6 = $ws.Range("D6").Style
$ws.Range("D6").Value = '''131.95'
$ws.Range("D6").Style = $styleD6
$ws.Range("E6").Value = '  +0.60%  '
$ws.Range("E7").Value = '  +0.02%  '
$styleD8 = $ws.Range("D8").Style
$ws.Range("D8").Value = '''0.587'
$ws.Range("D8").Style = $styleD8
$ws.Range("E8").Value = '  -0.19%  '
$styleD9 = $ws.Range("D9").Style
$ws.Range("D9").Value = '''2.331.80'
$ws.Range("D9").Style = $styleD9
$ws.Range("E9").Value = '  +2.10%  '
$ws.Range("E10").Value = '  +1.17%  '
$styleD11 = $ws.Range("D11").Style
$ws.Range("D11").Value = '''5.51'
$ws.Range("D11").Style = $styleD11
$ws.Range("E11").Value = '  +1.31%  '
$ws.Range("E12").Value = '  +0.90%  '
$styleD13 = $ws.Range("D13").Style
$ws.Range("D13").Value = '''0.333'
$ws.Range("D13").Style = $styleD13
$ws.Range("E13").Value = '  +0.89%  '
$styleD14 = $ws.Range("D14").Style
$ws.Range("D14").Value = '''23.83'
$ws.Range("D14").Style = $styleD14
$ws.Range("E14").Value = '  +1.30%  '
$styleD15 = $ws.Range("D15").Style
$ws.Range("D15").Value = '''2.754.60'
$ws.Range("D15").Style = $styleD15
$ws.Range("E15").Value = '  +2.24%  '
$styleD16 = $ws.Range("D16").Style
$ws.Range("D16").Value = '''59.903.82'
$ws.Range("D16").Style = $styleD16
$ws.Range("E16").Value = '  +3.43%  '
$ws.Range("E17").Value = '  +1.22%  '
$styleD18 = $ws.Range("D18").Style
$ws.Range("D18").Value = '''2.336.07'
$ws.Range("D18").Style = $styleD18
$ws.Range("E18").Value = '  +1.85%  '
$styleD19 = $ws.Range("D19").Style
$ws.Range("D19").Value = '''10.65'
$ws.Range("D19").Style = $styleD19
$ws.Range("E19").Value = '  +1.13%  '
$styleD20 = $ws.Range("D20").Style
$ws.Range("D20").Value = '''4.16'
$ws.Range("D20").Style = $styleD20
$ws.Range("E20").Value = '  -1.58%  '
$styleD21 = $ws.Range("D21").Style
$ws.Range("D21").Value = '''6.81'
$ws.Range("D21").Style = $styleD21
$ws.Range("E21").Value = '  +6.93%  '
$styleD22 = $ws.Range("D22").Style
$ws.Range("D22").Value = '''314.06'
$ws.Range("D22").Style = $styleD22
$ws.Range("E22").Value = '  +0.38%  '
$styleD23 = $ws.Range("D23").Style
$ws.Range("D23").Value = '''0.997'
$ws.Range("D23").Style = $styleD23
$ws.Range("E23").Value = '  -0.39%  '
$styleD24 = $ws.Range("D24").Style
$ws.Range("D24").Value = '''63.06'
$ws.Range("D24").Style = $styleD24
$ws.Range("E24").Value = '  +1.02%  '
$styleD25 = $ws.Range("D25").Style
$ws.Range("D25").Value = '''0.170'
$ws.Range("D25").Style = $styleD25
$ws.Range("E25").Value = '  +2.04%  '
$styleD26 = $ws.Range("D26").Style
$ws.Range("D26").Value = '''1.01'
$ws.Range("D26").Style = $styleD26
$ws.Range("E26").Value = '  +0.59%  '
$styleD27 = $ws.Range("D27").Style
$ws.Range("D27").Value = '''7.87'
$ws.Range("D27").Style = $styleD27
$ws.Range("E27").Value = '  -1.48%  '
$ws.Range("E28").Value = '  +4.74%  '
$styleD29 = $ws.Range("D29").Style
$ws.Range("D29").Value = '''1.74'
$ws.Range("D29").Style = $styleD29
$ws.Range("E29").Value = '  +2.25%  '
$styleD30 = $ws.Range("D30").Style
$ws.Range("D30").Value = '''171.17'
$ws.Range("D30").Style = $styleD30
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("E31").Value = '  +8.81%  '
$styleD32 = $ws.Range("D32").Style
$ws.Range("D32").Value = '''0.0₃0727'
$ws.Range("D32").Style = $styleD32
$ws.Range("E32").Value = '  +1.10%  '
$styleD33 = $ws.Range("D33").Style
$ws.Range("D33").Value = '''5.92'
$ws.Range("D33").Style = $styleD33
$ws.Range("E33").Value = '  +2.81%  '
$ws.Range("E34").Value = '  +14.38%  '
$styleD35 = $ws.Range("D35").Style
$ws.Range("D35").Value = '''0.380'
$ws.Range("D35").Style = $styleD35
$ws.Range("E35").Value = '  +0.46%  '
$styleD36 = $ws.Range("D36").Style
$ws.Range("D36").Value = '''17.99'
$ws.Range("D36").Style = $styleD36
$ws.Range("E36").Value = '  +1.42%  '
$ws.Range("E37").Value = '  +0.00%  '
$styleD38 = $ws.Range("D38").Style
$ws.Range("D38").Value = '''0.998'
$ws.Range("D38").Style = $styleD38
$ws.Range("E38").Value = '  -0.22%  '
$styleD39 = $ws.Range("D39").Style
$ws.Range("D39").Value = '''4.14'
$ws.Range("D39").Style = $styleD39
$ws.Range("E39").Value = '  +6.45%  '
$styleD40 = $ws.Range("D40").Style
$ws.Range("D40").Value = '''315.33'
$ws.Range("D40").Style = $styleD40
$ws.Range("E40").Value = '  +9.98%  '
$styleD41 = $ws.Range("D41").Style
$ws.Range("D41").Value = '''38.13'
$ws.Range("D41").Style = $styleD41
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("E42").Value = '  +2.77%  '
$styleD43 = $ws.Range("D43").Style
$ws.Range("D43").Value = '''142.45'
$ws.Range("D43").Style = $styleD43
$ws.Range("E43").Value = '  +0.44%  '
$styleD44 = $ws.Range("D44").Style
$ws.Range("D44").Value = '''3.45'
$ws.Range("D44").Style = $styleD44
$ws.Range("E44").Value = '  +0.59%  '
$styleD45 = $ws.Range("D45").Style
$ws.Range("D45").Value = '''0.0952'
$ws.Range("D45").Style = $styleD45
$ws.Range("E45").Value = '  +0.07%  '
$styleD46 = $ws.Range("D46").Style
$ws.Range("D46").Value = '''19.23'
$ws.Range("D46").Style = $styleD46
$ws.Range("E46").Value = '  +6.53%  '
$styleD47 = $ws.Range("D47").Style
$ws.Range("D47").Value = '''0.0496'
$ws.Range("D47").Style = $styleD47
$ws.Range("E47").Value = '  +0.26%  '
$styleD48 = $ws.Range("D48").Style
$ws.Range("D48").Value = '''0.558'
$ws.Range("D48").Style = $styleD48
$ws.Range("E48").Value = '  +1.04%  '
$ws.Range("E49").Value = '  +1.43%  '
$styleD50 = $ws.Range("D50").Style
$ws.Range("D50").Value = '''11.01'
$ws.Range("D50").Style = $styleD50
$ws.Range("E50").Value = '  +0.88%  '
$styleD51 = $ws.Range("D51").Style
$ws.Range("D51").Value = '''0.0₆0204'
$ws.Range("D51").Style = $styleD51
$ws.Range("E51").Value = '  +11.52%  '
